# Generate Report for Handoff
# Update the localization-status workbook to reflect that the
# 9b2c1858-178e-4bff-b5ee-50b671b78afc.md file is now "Ready for handoff"
# (previously "Handed back: in sync with en-US"), with updated handoff
# timestamps on the per-locale sheets and the rolled-up Overview sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 9b2c1858-... (row 3) ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-21 22:46:10"

# --- zh-cn sheet: row for 9b2c1858-... (row 3) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-21 22:46:03"

# --- de-de sheet: row for 9b2c1858-... (row 3) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-21 22:46:10"
